# Update column F (dSF) values on the active sheet to reflect the
# repulled / recalculated data, per the commit:
# "repull data, push all data, mean calculation"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 4
    5  = 5
    6  = -3
    8  = -5
    9  = -3
    11 = 0
    13 = -3
    14 = 5
    15 = -2
    16 = -4
    17 = -4
    18 = -2
    19 = -2
    20 = 2
    21 = -2
    22 = -2
    24 = -1
    25 = 0
    26 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
